$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text updates ---
$ws.Range("A8").Value = "Volume 31   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/8/2024  Through  1/14/2024"

# --- Style fixups for cells changing from text(NA) to number ---
$ws.Range("D22").NumberFormat = $ws.Range("D16").NumberFormat
$ws.Range("J22").NumberFormat = $ws.Range("J16").NumberFormat
$ws.Range("C26").NumberFormat = $ws.Range("C16").NumberFormat
$ws.Range("I26").NumberFormat = $ws.Range("I16").NumberFormat
$ws.Range("F30").NumberFormat = $ws.Range("F16").NumberFormat
$ws.Range("I30").NumberFormat = $ws.Range("I16").NumberFormat
$ws.Range("E22").NumberFormat = $ws.Range("E16").NumberFormat
$ws.Range("K22").NumberFormat = $ws.Range("K16").NumberFormat

# --- Numeric cell value updates ---
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 4
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -16.666666666666
$ws.Range("F16").Value = 21
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 9
$ws.Range("J16").Value = 11
$ws.Range("K16").Value = -18.181818181818
$ws.Range("L16").Value = 12.5
$ws.Range("M16").Value = -35.714285714285
$ws.Range("N16").Value = -83.018867924528
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 300
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 14.285714285714
$ws.Range("I17").Value = 9
$ws.Range("J17").Value = 8
$ws.Range("K17").Value = 12.5
$ws.Range("L17").Value = 125
$ws.Range("M17").Value = 80
$ws.Range("N17").Value = -30.76923076923
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 30.76923076923
$ws.Range("I18").Value = 10
$ws.Range("J18").Value = 9
$ws.Range("K18").Value = 11.111111111111
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 150
$ws.Range("N18").Value = -82.758620689655
$ws.Range("C19").Value = 19
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 46.153846153846
$ws.Range("F19").Value = 66
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = 20
$ws.Range("I19").Value = 35
$ws.Range("J19").Value = 25
$ws.Range("K19").Value = 40
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 150
$ws.Range("N19").Value = 12.903225806451
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 3
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = 5.882352941176
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = 5
$ws.Range("K20").Value = -20
$ws.Range("M20").Value = -20
$ws.Range("N20").Value = -95.402298850574
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = 13.793103448275
$ws.Range("F21").Value = 140
$ws.Range("G21").Value = 117
$ws.Range("H21").Value = 19.658119658119
$ws.Range("I21").Value = 67
$ws.Range("J21").Value = 62
$ws.Range("K21").Value = 8.064516129032
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 55.813953488372
$ws.Range("N21").Value = -72.427983539094
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = -40
$ws.Range("I22").Value = 2
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = -33.333333333333
$ws.Range("M22").Value = -33.333333333333
$ws.Range("C24").Value = 38
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = 46.153846153846
$ws.Range("F24").Value = 174
$ws.Range("G24").Value = 115
$ws.Range("H24").Value = 51.304347826087
$ws.Range("I24").Value = 74
$ws.Range("J24").Value = 47
$ws.Range("K24").Value = 57.446808510638
$ws.Range("L24").Value = 100
$ws.Range("M24").Value = 174.074074074074
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = -68.75
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = 10.81081081081
$ws.Range("I25").Value = 20
$ws.Range("J25").Value = 24
$ws.Range("K25").Value = -16.666666666666
$ws.Range("L25").Value = 17.647058823529
$ws.Range("M25").Value = 66.666666666666
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -50
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = -40
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = 5
$ws.Range("K26").Value = -80
$ws.Range("L26").Value = -50
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -42.857142857142
$ws.Range("I27").Value = 3
$ws.Range("J27").Value = 3
$ws.Range("L27").Value = 50
$ws.Range("F30").Value = 1
$ws.Range("I30").Value = 1
